# "A Lot of Attack Elements Done"
# Fill in the Get-Damage-Gif and Death-Gif addresses for the Siavash card,
# and leave the view scrolled/selected near those columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 previously had the Attack Gif path copy/pasted into the Get Damage
# Gif (O2) and Death Gif (P2) columns too - replace those placeholders with
# their real gif addresses.
$ws.Range("O2").Value = "./res/gifs/f1_altgeneral/hit_t.gif"
$ws.Range("P2").Value = "./res/gifs/f1_altgeneral/death_t.gif"

# Scroll the view over towards the newly-edited columns and leave the
# selection on N27, matching where editing left off.
$win = $excel.ActiveWindow
$win.ScrollColumn = 11
$win.ScrollRow = 1
$ws.Range("N27").Select()
